$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear contents of row 2 (23U2292 / Djine Sinto Pafing / 45 / Bien) - row stays empty, not shifted
$ws.Range("A2:E2").ClearContents()

# Clear contents of row 4 (23U2355 / Martial Jeannot Maa / 43 / Bien)
$ws.Range("A4:E4").ClearContents()

# Delete column E (Observation) entirely
$ws.Columns.Item(5).Delete()

$ws.Columns.Item(5).Select()
